$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 426, shifting existing rows 426:511 down to 427:512.
$ws.Rows(426).Insert()

# Populate the newly inserted row 426 with the new data record.
$ws.Cells.Item(426, 1).Value = 6
$ws.Cells.Item(426, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(426, 3).Value = "Metropolitana"
$ws.Cells.Item(426, 4).Value = 45275
$ws.Cells.Item(426, 5).Value = 13
$ws.Cells.Item(426, 6).Value = 100112026
$ws.Cells.Item(426, 7).Value = "Haba"
$ws.Cells.Item(426, 8).Value = "Sin especificar"
$ws.Cells.Item(426, 9).Value = "Primera"
$ws.Cells.Item(426, 10).Value = 480
$ws.Cells.Item(426, 11).Value = 14000
$ws.Cells.Item(426, 12).Value = 15000
$ws.Cells.Item(426, 13).Value = 14521
$ws.Cells.Item(426, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(426, 15).Value = "Carahue"
$ws.Cells.Item(426, 16).Value = 581
$ws.Cells.Item(426, 17).Value = 25
$ws.Cells.Item(426, 18).Value = "Hortaliza"
